$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TELEXo")

# Update quantities in rows 24 and 25
$ws.Range("A24").Value = 8
$ws.Range("A25").Value = 8

# Move the "R102,R103,R104,R105" designator group from G25 (TopBottom) up to G24
$ws.Range("G24").Value = $ws.Range("G25").Value()
$ws.Range("G25").Clear()

# Update view state: scroll so row 2 is at top, and select A25
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("A25").Select()
